$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Each entry: row, year, new "Data" value (as text), isNewRow flag
$updates = @(
    @(2, 1800, "853", 0),
    @(12, 1810, "853", 0),
    @(13, 1811, "815", 0),
    @(14, 1812, "802", 0),
    @(15, 1813, "808", 0),
    @(16, 1814, "813", 0),
    @(17, 1815, "816", 0),
    @(18, 1816, "789", 0),
    @(19, 1817, "781", 0),
    @(20, 1818, "822", 0),
    @(21, 1819, "787", 0),
    @(22, 1820, "824", 0),
    @(23, 1821, "805", 0),
    @(24, 1822, "830", 0),
    @(25, 1823, "797", 0),
    @(26, 1824, "786", 0),
    @(27, 1825, "853", 0),
    @(28, 1826, "885", 0),
    @(29, 1827, "904", 0),
    @(30, 1828, "910", 0),
    @(31, 1829, "932", 0),
    @(32, 1830, "925", 0),
    @(33, 1831, "929", 0),
    @(34, 1832, "956", 0),
    @(35, 1833, "961", 0),
    @(36, 1834, "983", 0),
    @(37, 1835, "1003", 0),
    @(38, 1836, "991", 0),
    @(39, 1837, "1028", 0),
    @(40, 1838, "1033", 0),
    @(41, 1839, "1023", 0),
    @(42, 1840, "1081", 0),
    @(43, 1841, "1081", 0),
    @(44, 1842, "1089", 0),
    @(45, 1843, "1111", 0),
    @(46, 1844, "1117", 0),
    @(47, 1845, "1141", 0),
    @(48, 1846, "1170", 0),
    @(49, 1847, "1183", 0),
    @(50, 1848, "1227", 0),
    @(51, 1849, "1298", 0),
    @(52, 1850, "1352", 0),
    @(53, 1851, "1392", 0),
    @(54, 1852, "1404", 0),
    @(55, 1853, "1364", 0),
    @(56, 1854, "1387", 0),
    @(57, 1855, "1425", 0),
    @(58, 1856, "1427", 0),
    @(59, 1857, "1462", 0),
    @(60, 1858, "1513", 0),
    @(61, 1859, "1559", 0),
    @(62, 1860, "1588", 0),
    @(63, 1861, "1575", 0),
    @(64, 1862, "1553", 0),
    @(65, 1863, "1596", 0),
    @(66, 1864, "1667", 0),
    @(67, 1865, "1710", 0),
    @(68, 1866, "1750", 0),
    @(69, 1867, "1666", 0),
    @(70, 1868, "1714", 0),
    @(71, 1869, "1849", 0),
    @(72, 1870, "1868", 0),
    @(73, 1871, "1836", 0),
    @(74, 1872, "1910", 0),
    @(75, 1873, "2024", 0),
    @(76, 1874, "1954", 0),
    @(77, 1875, "2050", 0),
    @(78, 1876, "1988", 0),
    @(79, 1877, "1908", 0),
    @(80, 1878, "1953", 0),
    @(81, 1879, "2219", 0),
    @(82, 1880, "2418", 0),
    @(83, 1881, "2488", 0),
    @(84, 1882, "2694", 0),
    @(85, 1883, "2683", 0),
    @(86, 1884, "2703", 0),
    @(87, 1885, "2547", 0),
    @(88, 1886, "2625", 0),
    @(89, 1887, "2783", 0),
    @(90, 1888, "2617", 0),
    @(91, 1889, "2708", 0),
    @(92, 1890, "2907", 0),
    @(93, 1891, "3038", 0),
    @(94, 1892, "3005", 0),
    @(95, 1893, "3084", 0),
    @(96, 1894, "2978", 0),
    @(97, 1895, "3304", 0),
    @(98, 1896, "3379", 0),
    @(99, 1897, "3210", 0),
    @(100, 1898, "3609", 0),
    @(101, 1899, "3467", 0),
    @(102, 1900, "3386", 0),
    @(103, 1901, "3416", 0),
    @(104, 1902, "3604", 0),
    @(105, 1903, "3381", 0),
    @(106, 1904, "3647", 0),
    @(107, 1905, "3470", 0),
    @(108, 1906, "3668", 0),
    @(109, 1907, "3810", 0),
    @(110, 1908, "4020", 0),
    @(111, 1909, "4042", 0),
    @(112, 1910, "4485", 0),
    @(113, 1911, "4481", 0),
    @(114, 1912, "4806", 0),
    @(115, 1913, "4836", 0),
    @(116, 1914, "4071", 0),
    @(117, 1915, "3813", 0),
    @(118, 1916, "4600", 0),
    @(119, 1917, "4621", 0),
    @(120, 1918, "4602", 0),
    @(121, 1919, "3853", 0),
    @(122, 1920, "4248", 0),
    @(123, 1921, "3727", 0),
    @(124, 1922, "3814", 0),
    @(125, 1923, "4514", 0),
    @(126, 1924, "4779", 0),
    @(127, 1925, "4903", 0),
    @(128, 1926, "4702", 0),
    @(129, 1927, "4544", 0),
    @(130, 1928, "5609", 0),
    @(131, 1929, "5679", 0),
    @(132, 1930, "4876", 0),
    @(133, 1931, "3904", 0),
    @(134, 1932, "2950", 0),
    @(135, 1933, "3618", 0),
    @(136, 1934, "4269", 0),
    @(137, 1935, "4511", 0),
    @(138, 1936, "4584", 0),
    @(139, 1937, "5083", 0),
    @(140, 1938, "5070", 0),
    @(141, 1939, "5106", 0),
    @(142, 1940, "5177", 0),
    @(143, 1941, "5091", 0),
    @(144, 1942, "5161", 0),
    @(145, 1943, "5209", 0),
    @(146, 1944, "5209", 0),
    @(147, 1945, "5553", 0),
    @(148, 1946, "5918", 0),
    @(149, 1947, "5187", 0),
    @(150, 1948, "5939", 0),
    @(151, 1949, "5710", 0),
    @(152, 1950, "5880", 0),
    @(153, 1951, "6001", 0),
    @(154, 1952, "6252", 0),
    @(155, 1953, "6582", 0),
    @(156, 1954, "6237", 0),
    @(157, 1955, "6341", 0),
    @(158, 1956, "6285", 0),
    @(159, 1957, "6758", 0),
    @(160, 1958, "6958", 0),
    @(161, 1959, "6409", 0),
    @(162, 1960, "6781", 0),
    @(163, 1961, "6923", 0),
    @(164, 1962, "7071", 0),
    @(165, 1963, "7336", 0),
    @(166, 1964, "7320", 0),
    @(167, 1965, "7208", 0),
    @(168, 1966, "7844", 0),
    @(169, 1967, "7933", 0),
    @(170, 1968, "8053", 0),
    @(171, 1969, "8187", 0),
    @(172, 1970, "8195", 0),
    @(173, 1971, "8773", 0),
    @(174, 1972, "8520", 0),
    @(175, 1973, "7911", 0),
    @(176, 1974, "7857", 0),
    @(177, 1975, "6731", 0),
    @(178, 1976, "6868", 0),
    @(179, 1977, "7438", 0),
    @(180, 1978, "7936", 0),
    @(181, 1979, "8475", 0),
    @(182, 1980, "9024", 0),
    @(183, 1981, "9427", 0),
    @(184, 1982, "8016", 0),
    @(185, 1983, "7667", 0),
    @(186, 1984, "7992", 0),
    @(187, 1985, "8024", 0),
    @(188, 1986, "8325", 0),
    @(189, 1987, "8721", 0),
    @(190, 1988, "9199", 0),
    @(191, 1989, "10005", 0),
    @(192, 1990, "10203", 0),
    @(193, 1991, "10746.4915641559", 0),
    @(194, 1992, "11773.3007091455", 0),
    @(195, 1993, "12296.9256463059", 0),
    @(196, 1994, "12692.5511321862", 0),
    @(197, 1995, "13715.9329954789", 0),
    @(198, 1996, "14433.3197742387", 0),
    @(199, 1997, "15079.8311950225", 0),
    @(200, 1998, "15259.5408475535", 0),
    @(201, 1999, "14846.4175283184", 0),
    @(202, 2000, "15211.6232956044", 0),
    @(203, 2001, "15448.008813318", 0),
    @(204, 2002, "15509.8247751724", 0),
    @(205, 2003, "15809.1777684547", 0),
    @(206, 2004, "16478.9727867052", 0),
    @(207, 2005, "17137.4860414664", 0),
    @(208, 2006, "17920.8031312197", 0),
    @(209, 2007, "18487.9007351472", 0),
    @(210, 2008, "18804.0293633295", 0),
    @(211, 2009, "18184.4813776621", 0),
    @(212, 2010, "18909.8472762053", 0),
    @(213, 2011, "19705", 1),
    @(214, 2012, "20531", 1),
    @(215, 2013, "21135", 1),
    @(216, 2014, "21335", 1),
    @(217, 2015, "21589", 1),
    @(218, 2016, "21696", 1)
)

foreach ($u in $updates) {
    $row = $u[0]
    $year = $u[1]
    $val = $u[2]
    $isNew = $u[3]
    if ($isNew -eq 1) {
        $ws.Cells.Item($row, 1).Value = 152
        $ws.Cells.Item($row, 2).Value = "Chile"
        $ws.Cells.Item($row, 3).Value = "GDP per Capita"
        $ws.Cells.Item($row, 4).Value = $year
    }
    $cell = $ws.Cells.Item($row, 5)
    # Force the numeric-looking value to be stored as text (matches the
    # "Data" column convention used throughout this sheet), then drop the
    # temporary formatting so no visible number format is left behind.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

